$wb = $excel.ActiveWorkbook

# Update info sheet B2 (timing value)
$wsInfo = $wb.Worksheets.Item("info")
$wsInfo.Cells.Item(2, 2).Value = [double]"2.759763717651367"

# Update x sheet column B (rows 2 to 151) with new initial-variable values
$wsX = $wb.Worksheets.Item("x")
$xValues = @(
    "-2.167974759139391",
    "-4.45116399200965e-08",
    "-4.415042258835782e-08",
    "-4.32354106087106e-08",
    "8.093617738804096e-08",
    "7.552474221313034e-09",
    "-3.280627424941995e-07",
    "-2.069951930614636e-08",
    "4.864829878640674e-08",
    "-4.959913618552829e-08",
    "-1.746030840521543e-07",
    "-1.833252864424049e-08",
    "9.915642948839831e-08",
    "2.648559854494112e-08",
    "-1.217567644560372e-08",
    "-6.606339109336398e-09",
    "6.339272923198112e-08",
    "-2.205945713536315e-07",
    "-6.691084765124262e-10",
    "-2.289464003823771e-08",
    "1.240619212193401e-07",
    "1.033543362961711e-08",
    "-1.713009658489686e-10",
    "4.340179400893044e-09",
    "-2.984389887770847e-08",
    "-2.772059944779669e-08",
    "-2.552036657816202e-07",
    "5.810589971611347e-09",
    "3.381507815914042e-08",
    "1.776441693050156e-07",
    "-1.517849802636626e-09",
    "-3.880955507940747e-08",
    "-1.714342229091974e-08",
    "-3.336034412824929e-08",
    "-3.519200925135336e-08",
    "-1.652879845655146e-08",
    "-1.143463536867534e-08",
    "3.700758963039585e-09",
    "5.936741019438803e-09",
    "-2.164249976821129e-08",
    "1.63100775687885e-07",
    "5.621459991639807e-12",
    "7.950757664511333e-11",
    "-5.735919433450573e-09",
    "-1.458574304283953e-08",
    "-1.055091490101525e-07",
    "-1.697561509188183e-07",
    "4.133592088596386e-08",
    "1.730007636702715e-08",
    "4.068535927044412e-08",
    "6.537085001944038e-09",
    "-1.15348482637547e-07",
    "1.001711075715983e-07",
    "3.709913978792273e-09",
    "-1.176679745659932e-08",
    "2.224995863497323e-08",
    "-1.270585541431134e-07",
    "8.106246251304943e-08",
    "-5.467600692097412e-09",
    "2.77386120893249e-08",
    "-9.186648532055886e-09",
    "4.037549136183775e-08",
    "-4.630292938992853e-08",
    "1.040243834863739e-07",
    "1.80645184861014e-09",
    "-1.077600528698684e-07",
    "-5.820828488543236e-08",
    "-1.808805226306442e-08",
    "1.024334378698204e-08",
    "-2.024487655249814e-07",
    "2.551681388648721e-08",
    "-3.928242378330924e-08",
    "2.272350492237694e-08",
    "7.235987430101805e-08",
    "-1.916089047968429e-07",
    "1.342259756741415e-09",
    "-9.497211131657001e-08",
    "1.722868296798226e-07",
    "-1.220171065940722e-07",
    "1.97935093018527e-08",
    "9.450571823418436e-08",
    "-7.871571021997854e-10",
    "1.440345818858124e-07",
    "-2.302148006456259",
    "1.092527592872604e-07",
    "5.371973091507048e-08",
    "5.399321283907277e-08",
    "-1.057833161400604e-05",
    "-7.060306122231917e-08",
    "-1.927809880499868e-08",
    "-3.588031859049051e-08",
    "5.096936872694205e-07",
    "9.182938726874839e-09",
    "3.180361619678821e-08",
    "5.487694959434483e-08",
    "2.312016401507387e-09",
    "-1.55395556723788e-08",
    "-2.896988541214729e-08",
    "3.753723692561368e-09",
    "-1.192961469286118e-07",
    "1.554995620528322e-09",
    "-1.394482328826442e-07",
    "4.932626708828706e-08",
    "5.976186072915335e-09",
    "1.598854102732089e-08",
    "6.53008529702027e-09",
    "-1.501150255591833e-09",
    "5.818087132873942e-08",
    "-1.750740436900304e-08",
    "2.278069960869292e-08",
    "-3.907208249400777e-08",
    "8.49082967928591e-08",
    "-2.591590937471511e-08",
    "-4.00330616844083e-08",
    "1.647403877851825e-09",
    "1.375301304799698e-06",
    "-8.916209504041733e-08",
    "-5.601657262944973e-09",
    "-2.145201800178086e-06",
    "4.722950850456638e-08",
    "-3.94357428398819e-08",
    "-3.413554619602102e-08",
    "2.904103326016521e-08",
    "-1.525959107995169e-07",
    "2.694446636031043e-08",
    "-6.309130831898615e-10",
    "-2.338792415000123e-07",
    "1.86536544625206e-08",
    "-1.624981930455002e-08",
    "-9.996524594646756e-09",
    "1.354523578124798e-07",
    "-1.31615459582131e-08",
    "-1.544385619589896e-07",
    "-5.836416983899436e-09",
    "-2.345553402253957e-06",
    "-2.063984807960972e-08",
    "-1.927689995805912e-07",
    "-5.5411686134191e-08",
    "-3.579383455385554e-07",
    "-1.750352539996732e-09",
    "-7.43242433491053e-08",
    "1.221113493756351e-08",
    "4.177898498185276e-08",
    "-3.114706906834342e-08",
    "4.831846846102918e-08",
    "-6.789960101107514e-09",
    "-2.913311648141036e-09",
    "1.810239916731048e-10",
    "-2.291411444444226e-07",
    "3.096566390328359e-09"
)

for ($i = 0; $i -lt $xValues.Count; $i++) {
    $row = $i + 2
    $wsX.Cells.Item($row, 2).Value = [double]$xValues[$i]
}

Write-Output "Updated info!B2 and x!B2:B151"
